$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Merk: PMI"
$ws.Range("H3").Value = "No.Lot: 71616155"
$ws.Range("P3").Value = "Exp: 12 desember 2026"

$ws.Range("M5").Select()
